# Applies the "Pipeline" table addition and minor view/formatting changes
# described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Pipeline" table in columns J:O, rows 1-14 ---

# Header row (bold, centered) - style "1" in the original file
$ws.Range("J1").Value = "Pipeline"
$ws.Range("K1").Value = "RF Read"
$ws.Range("L1").Value = "RF Write"
$ws.Range("M1").Value = "RAM Read"
$ws.Range("N1").Value = "RAM Write"
$ws.Range("O1").Value = "PC Write"
$ws.Range("J1:O1").Font.Bold = $true
$ws.Range("J1:O1").HorizontalAlignment = -4108  # xlCenter

# Data rows
$ws.Range("J2").Value = "NOP"

$ws.Range("J3").Value = "ALU"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1

$ws.Range("J4").Value = "RDMi"
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1

$ws.Range("J5").Value = "WRMi"
$ws.Range("K5").Value = 1
$ws.Range("N5").Value = 1

$ws.Range("J6").Value = "IML"
$ws.Range("L6").Value = 1

$ws.Range("J7").Value = "IMH"
$ws.Range("L7").Value = 1

$ws.Range("J8").Value = "RDMr"
$ws.Range("K8").Value = 1
$ws.Range("M8").Value = 1

$ws.Range("J9").Value = "WRMr"
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = 1

$ws.Range("J10").Value = "JC"
$ws.Range("K10").Value = 1
$ws.Range("O10").Value = 1

$ws.Range("J11").Value = "JR"
$ws.Range("K11").Value = 1
$ws.Range("O11").Value = 1

$ws.Range("J12").Value = "JA"
$ws.Range("O12").Value = 1

$ws.Range("J13").Value = "CR"
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1

$ws.Range("J14").Value = "EXT"
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 1

# --- View / formatting tweaks ---
# Scroll the window so column F becomes the leftmost visible column,
# zoom to 145%, and leave the selection on N11 (matches the saved view).
$excel.ActiveWindow.Zoom = 145
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N11").Select() | Out-Null
